$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (current Seed column) to hold "sec/run"
$ws.Columns("C:C").Insert()

# New header for column C
$ws.Range("C1").Value = "sec/run"

# Formulas for sec/run = Time / Runs
$ws.Range("C2").Formula = "=B2/A2"

# Row 3: B3 is a string ("?"), so no formula is placed in C3 for this row (left blank)

# Existing data shifted: D3 previously held 595836831489578 as C3 (pre-insert); now update values that differ from a pure shift
$ws.Range("D3").Value = 5065
$ws.Range("E3").Value = 1

# New rows 4 and 5
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = 1123.7
$ws.Range("C4").Formula = "=B4/A4"
$ws.Range("D4").Value = 1870500046
$ws.Range("E4").Value = 16

$ws.Range("A5").Value = 20
$ws.Range("C5").Formula = "=B5/A5"
$ws.Range("D5").Value = 1448491702
$ws.Range("E5").Value = 31

# sdlist header moved from E3 to F3 automatically by column insert; make sure F3 has the sdlist string still
# (should already be there via shift, but set explicitly to be safe)
$ws.Range("F3").Value = "sdlist"

# Column width for column D (Seed) per diff: width 11, bestFit
# (the engine stores widths with the standard ~0.83 padding baked in, so we
# dial in the COM-level ColumnWidth that round-trips to an OOXML width of 11)
$ws.Columns("D:D").ColumnWidth = 10.1666666666667

# Selection per diff
$ws.Range("D11").Select()
